$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): K1..O1 get reshuffled so "percentage"/"multiplier"
#     (previously N1/O1) move to K1/L1, and image name / library flag /
#     library_base_price (previously K1/L1/M1) shift right to M1/N1/O1.
$row1_K = $ws.Cells.Item(1, 11).Value2   # image name
$row1_L = $ws.Cells.Item(1, 12).Value2   # library (0 = no, 1 = yes)
$row1_M = $ws.Cells.Item(1, 13).Value2   # library_base_price
$row1_N = $ws.Cells.Item(1, 14).Value2   # percentage
$row1_O = $ws.Cells.Item(1, 15).Value2   # multiplier

$ws.Cells.Item(1, 11).Value2 = $row1_N   # K1 = percentage
$ws.Cells.Item(1, 12).Value2 = $row1_O   # L1 = multiplier
$ws.Cells.Item(1, 13).Value2 = $row1_K   # M1 = image name
$ws.Cells.Item(1, 14).Value2 = $row1_L   # N1 = library (0 = no, 1 = yes)
$ws.Cells.Item(1, 15).Value2 = $row1_M   # O1 = library_base_price

# --- Row 2: K2/L2 (image name, library flag) shift to M2/N2.
$row2_K = $ws.Cells.Item(2, 11).Value2
$row2_L = $ws.Cells.Item(2, 12).Value2

$ws.Cells.Item(2, 11).ClearContents()
$ws.Cells.Item(2, 12).ClearContents()
$ws.Cells.Item(2, 13).Value2 = $row2_K
$ws.Cells.Item(2, 14).Value2 = $row2_L

# --- Row 3: K3/L3/M3 shift to M3/N3/O3; the old N3/O3 (percentage,
#     multiplier sample values) are dropped entirely.
$row3_K = $ws.Cells.Item(3, 11).Value2
$row3_L = $ws.Cells.Item(3, 12).Value2
$row3_M = $ws.Cells.Item(3, 13).Value2

$ws.Cells.Item(3, 11).ClearContents()
$ws.Cells.Item(3, 12).ClearContents()
$ws.Cells.Item(3, 13).Value2 = $row3_K
$ws.Cells.Item(3, 14).Value2 = $row3_L   # overwrites old N3 (0.8) sample
$ws.Cells.Item(3, 15).Value2 = $row3_M   # overwrites old O3 (2) sample

# --- Column widths: new bestFit columns for A, K, L; O keeps the old
#     13-15 width, now scoped to just column O.
$ws.Columns.Item(1).ColumnWidth = 8.944010416666666    # -> ~9.77734375
$ws.Columns.Item(11).ColumnWidth = 9.276041666666666   # -> ~10.109375
$ws.Columns.Item(12).ColumnWidth = 7.830729166666667   # -> ~8.6640625
$ws.Columns.Item(15).ColumnWidth = 13.833333333333334  # -> ~14.6640625

# --- View: scroll so column M is the leftmost visible column, and
#     select the header row A1:S1.
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
$ws.Range("A1:S1").Select()
